$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the IFRS financial data block (rows 2-9, columns D:AJ).
# Rows 2-6: replace the mis-scaled figures with the corrected values
#           (a few metrics that are no longer available are cleared).
# Rows 7-9: the whole data block (D:AJ) is cleared, leaving only the
#           row label columns A:C, since those periods have no data now.

# --- Row 2 ---
$ws.Range("D2").Value = 241
$ws.Range("E2").Value = -91
$ws.Range("F2").Value = -91
$ws.Range("G2").Value = -130
$ws.Range("H2").Value = -130
$ws.Range("I2").Value = -130
$ws.Range("K2").Value = 648
$ws.Range("L2").Value = 269
$ws.Range("M2").Value = 379
$ws.Range("N2").Value = 379
$ws.Range("P2").Value = 445
$ws.Range("Q2").Value = -62
$ws.Range("R2").Value = -19
$ws.Range("S2").Value = 74
$ws.Range("T2").Value = 7
$ws.Range("U2").Value = -69
$ws.Range("V2").Value = 112
$ws.Range("W2").Value = -37.64
$ws.Range("X2").Value = -54.01
$ws.Range("Y2").Value = -38.04
$ws.Range("Z2").Value = -18.5
$ws.Range("AA2").Value = 71.01000000000001
$ws.Range("AB2").Value = -7.19
$ws.Range("AC2").Value = -30
$ws.Range("AD2").Value = -4.52
$ws.Range("AE2").Value = 75
$ws.Range("AF2").Value = 1.81
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 502923091
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# --- Row 3 ---
$ws.Range("D3").Value = 255
$ws.Range("E3").Value = -24
$ws.Range("F3").Value = -24
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 22
$ws.Range("I3").Value = 21
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 701
$ws.Range("L3").Value = 219
$ws.Range("M3").Value = 483
$ws.Range("N3").Value = 482
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 501
$ws.Range("Q3").Value = -84
$ws.Range("R3").Value = 29
$ws.Range("S3").Value = 78
$ws.Range("T3").Value = 9
$ws.Range("U3").Value = -94
$ws.Range("V3").Value = 110
$ws.Range("W3").Value = -9.57
$ws.Range("X3").Value = 8.43
$ws.Range("Y3").Value = 4.94
$ws.Range("Z3").Value = 3.19
$ws.Range("AA3").Value = 45.3
$ws.Range("AB3").Value = 2.25
$ws.Range("AC3").Value = 4
$ws.Range("AD3").Value = 100.47
$ws.Range("AE3").Value = 85
$ws.Range("AF3").Value = 4.92
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 566138828

# --- Row 4 ---
$ws.Range("D4").Value = 233
$ws.Range("E4").Value = -86
$ws.Range("F4").Value = -86
$ws.Range("G4").Value = -95
$ws.Range("H4").Value = -95
$ws.Range("I4").Value = -95
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 845
$ws.Range("L4").Value = 114
$ws.Range("M4").Value = 731
$ws.Range("N4").Value = 731
$ws.Range("P4").Value = 631
$ws.Range("Q4").Value = -86
$ws.Range("R4").Value = -42
$ws.Range("S4").Value = 237
$ws.Range("T4").Value = 5
$ws.Range("U4").Value = -91
$ws.Range("V4").Value = 5
$ws.Range("W4").Value = -36.71
$ws.Range("X4").Value = -40.71
$ws.Range("Y4").Value = -15.66
$ws.Range("Z4").Value = -12.27
$ws.Range("AA4").Value = 15.61
$ws.Range("AB4").Value = 20.54
$ws.Range("AC4").Value = -16
$ws.Range("AD4").Value = -17.05
$ws.Range("AE4").Value = 108
$ws.Range("AF4").Value = 2.48
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 679955804
$ws.Range("O4").ClearContents()

# --- Row 5 ---
$ws.Range("D5").Value = 345
$ws.Range("E5").Value = -95
$ws.Range("F5").Value = -95
$ws.Range("G5").Value = -137
$ws.Range("H5").Value = -137
$ws.Range("I5").Value = -137
$ws.Range("K5").Value = 845
$ws.Range("L5").Value = 238
$ws.Range("M5").Value = 606
$ws.Range("N5").Value = 606
$ws.Range("P5").Value = 631
$ws.Range("Q5").Value = -9
$ws.Range("R5").Value = -172
$ws.Range("S5").Value = 95
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = -14
$ws.Range("V5").Value = 76
$ws.Range("W5").Value = -27.44
$ws.Range("X5").Value = -39.61
$ws.Range("Y5").Value = -20.43
$ws.Range("Z5").Value = -16.17
$ws.Range("AA5").Value = 39.33
$ws.Range("AB5").Value = -0.74
$ws.Range("AC5").Value = -20
$ws.Range("AD5").Value = -10.27
$ws.Range("AE5").Value = 89
$ws.Range("AF5").Value = 2.31
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 679955804
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# --- Row 6 ---
$ws.Range("D6").Value = 345
$ws.Range("E6").Value = -114
$ws.Range("F6").Value = -114
$ws.Range("G6").Value = -133
$ws.Range("H6").Value = -133
$ws.Range("I6").Value = -133
$ws.Range("K6").Value = 687
$ws.Range("L6").Value = 211
$ws.Range("M6").Value = 476
$ws.Range("N6").Value = 476
$ws.Range("P6").Value = 635
$ws.Range("Q6").Value = -15
$ws.Range("R6").Value = -19
$ws.Range("S6").Value = 7
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = -17
$ws.Range("V6").Value = 89
$ws.Range("W6").Value = -33.16
$ws.Range("X6").Value = -38.56
$ws.Range("Y6").Value = -24.58
$ws.Range("Z6").Value = -17.38
$ws.Range("AA6").Value = 44.23
$ws.Range("AB6").Value = -12.75
$ws.Range("AC6").Value = -20
$ws.Range("AD6").Value = -8.59
$ws.Range("AE6").Value = 70
$ws.Range("AF6").Value = 2.41
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 684267382
$ws.Range("AI6").ClearContents()

# --- Rows 7-9: clear the whole data block, keep only A:C labels ---
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
